$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text price strings (some look numeric, e.g. "211.96").
# Force the whole price column to Text before writing so Excel doesn't
# auto-convert single-dot values into numbers; restore the default
# (unstyled) cell style afterwards so no stray formatting diff remains.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.259.62"
$ws.Range("E2").Value = "  +1.21%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.618.06"
$ws.Range("E3").Value = "  +1.93%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "211.96"
$ws.Range("E5").Value = "  +0.83%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.03%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.486"
$ws.Range("E7").Value = "  +0.86%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +0.95%  "

# Row 10 - Solana
$ws.Range("D10").Value = "18.75"
$ws.Range("E10").Value = "  +5.29%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0815"
$ws.Range("E11").Value = "  +0.83%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.843.34"

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.628.93"
$ws.Range("E13").Value = "  +2.50%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "3.99"
$ws.Range("E14").Value = "  +0.61%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "26.269.18"

# Row 17 - Litecoin
$ws.Range("D17").Value = "62.33"
$ws.Range("E17").Value = "  +4.08%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.0₃0726"
$ws.Range("E18").Value = "  +1.09%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.02%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "201.15"
$ws.Range("E20").Value = "  +1.59%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +1.74%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "9.30"
$ws.Range("E22").Value = "  +1.31%  "

# Row 23 - Chainlink
$ws.Range("E23").Value = "  +1.32%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +3.16%  "

# Row 25 - Monero
$ws.Range("D25").Value = "144.06"
$ws.Range("E25").Value = "  +0.70%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.07%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -1.12%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  +0.84%  "

# Row 29 - Cosmos
$ws.Range("D29").Value = "6.53"
$ws.Range("E29").Value = "  +1.78%  "

# Row 30 - Hedera
$ws.Range("D30").Value = "0.0516"
$ws.Range("E30").Value = "  +8.97%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.82%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +2.02%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -0.38%  "

# Row 34 - now HuobiToken (was LidoDAOToken)
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "2.41"
$ws.Range("E34").Value = "  +1.90%  "

# Row 35 - now LidoDAOToken (was HuobiToken)
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "1.49"
$ws.Range("E35").Value = "  +1.84%  "

# Row 36 - Maker
$ws.Range("D36").Value = "1.177.22"
$ws.Range("E36").Value = "  +4.92%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "0.0163"
$ws.Range("E37").Value = "  +1.24%  "

# Row 38 - ARBITRUM
$ws.Range("E38").Value = "  +3.17%  "

# Row 40 - MXToken
$ws.Range("E40").Value = "  +0.17%  "

# Row 41 - ImmutableX
$ws.Range("D41").Value = "0.494"
$ws.Range("E41").Value = "  +1.38%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "0.792"
$ws.Range("E42").Value = "  +1.27%  "

# Row 43 - FraxShare
$ws.Range("D43").Value = "5.33"
$ws.Range("E43").Value = "  +5.08%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.755.57"
$ws.Range("E44").Value = "  +1.85%  "

# Row 45 - Quant
$ws.Range("D45").Value = "92.73"
$ws.Range("E45").Value = "  +1.11%  "

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = "  +13.46%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  +2.36%  "

# Row 48 - Aave
$ws.Range("D48").Value = "53.67"
$ws.Range("E48").Value = "  +1.12%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  +1.08%  "

# Row 50 - Mantle
$ws.Range("D50").Value = "0.408"
$ws.Range("E50").Value = "  +0.40%  "

# Row 51 - USDD
$ws.Range("E51").Value = "  -0.23%  "

# Restore the default (unstyled) style on the whole price column now that
# every text value has been committed, so formatting matches the original.
$ws.Range("D2:D51").Style = "Normal"
